$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '69.832.17'
$r.Style = "Normal"
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  -1.93%  '
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '3.575.88'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  -2.47%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  -0.06%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '575.78'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  -3.41%  '
$r.Style = "Normal"
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '186.88'
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  -4.51%  '
$r.Style = "Normal"
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '3.569.15'
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  -2.53%  '
$r.Style = "Normal"
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -4.19%  '
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  +0.08%  '
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  +1.34%  '
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  -3.93%  '
$r.Style = "Normal"
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '55.27'
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  -6.51%  '
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  +3.64%  '
$r.Style = "Normal"
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '9.58'
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  -4.50%  '
$r.Style = "Normal"
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '4.153.72'
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  -2.40%  '
$r.Style = "Normal"
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  -1.80%  '
$r.Style = "Normal"
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '3.576.97'
$r.Style = "Normal"
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  -2.62%  '
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '69.868.99'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  -1.93%  '
$r.Style = "Normal"
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '12.55'
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  -2.29%  '
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  -0.58%  '
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -3.94%  '
$r.Style = "Normal"
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '489.40'
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  -1.04%  '
$r.Style = "Normal"
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '19.53'
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '4.93'
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  -8.70%  '
$r.Style = "Normal"
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '95.93'
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  +4.27%  '
$r.Style = "Normal"
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '4.32'
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  -5.00%  '
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  -7.08%  '
$r.Style = "Normal"
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '11.03'
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  -4.73%  '
$r.Style = "Normal"
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '9.28'
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  -4.28%  '
$r.Style = "Normal"
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '31.82'
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  -4.13%  '
$r.Style = "Normal"
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '7.54'
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  -5.18%  '
$r.Style = "Normal"
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '66.86'
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  +1.28%  '
$r.Style = "Normal"
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '12.03'
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -2.56%  '
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  -4.90%  '
$r.Style = "Normal"
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '570.50'
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -9.97%  '
$r.Style = "Normal"
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '38.25'
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  -6.40%  '
$r.Style = "Normal"
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  -0.02%  '
$r.Style = "Normal"
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '3.06'
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  +6.44%  '
$r.Style = "Normal"
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '0.0₃0799'
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  -5.26%  '
$r.Style = "Normal"
$r = $ws.Range('B40')
$r.NumberFormat = "@"
$r.Value = 'TheGraph'
$r.Style = "Normal"
$r = $ws.Range('C40')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '0.393'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  -5.49%  '
$r.Style = "Normal"
$r = $ws.Range('B41')
$r.NumberFormat = "@"
$r.Value = 'dogwifhat'
$r.Style = "Normal"
$r = $ws.Range('C41')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$r.Style = "Normal"
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '3.29'
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  +12.82%  '
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  -2.26%  '
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  -8.17%  '
$r.Style = "Normal"
$r = $ws.Range('B44')
$r.NumberFormat = "@"
$r.Value = 'Maker'
$r.Style = "Normal"
$r = $ws.Range('C44')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$r.Style = "Normal"
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '3.247.79'
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -2.53%  '
$r.Style = "Normal"
$r = $ws.Range('B45')
$r.NumberFormat = "@"
$r.Value = 'ThetaToken'
$r.Style = "Normal"
$r = $ws.Range('C45')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$r.Style = "Normal"
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '2.99'
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -5.82%  '
$r.Style = "Normal"
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '3.45'
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  +4.31%  '
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '9.68'
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  +3.68%  '
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  -4.47%  '
$r.Style = "Normal"
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '0.136'
$r.Style = "Normal"
$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  -2.20%  '
$r.Style = "Normal"
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  -0.22%  '
$r.Style = "Normal"
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '3.17'
$r.Style = "Normal"
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  -5.00%  '
$r.Style = "Normal"
